# Weekly refresh of the "Fruta, Terminal La Palmera de La Serena - Membrillo" data set.
# The market data (prices/volumes/origin/dates) for rows 2-20 is updated in place to
# reflect the latest weekly extract, and three additional records (rows 21-23) are
# appended, extending the used range from A1:T20 to A1:T23.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44630, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 16, 400000, 405000, 402500, '$/bins (450 kilos)', 'Región de O''Higgins', 894, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44309, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 20, 305000, 310000, 307500, '$/bins (450 kilos)', 'Provincia de Cachapoal', 683, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44309, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 20, 285000, 290000, 287500, '$/bins (450 kilos)', 'Provincia de Cachapoal', 639, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44309, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Segunda', 20, 255000, 260000, 257500, '$/bins (450 kilos)', 'Provincia de Cachapoal', 572, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44295, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 20, 225000, 230000, 227500, '$/bins (450 kilos)', 'Región Metropolitana', 506, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44295, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Segunda', 16, 195000, 200000, 197500, '$/bins (450 kilos)', 'Región Metropolitana', 439, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44655, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 20, 335000, 340000, 337500, '$/bins (450 kilos)', 'Región de O''Higgins', 750, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44655, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 20, 295000, 300000, 297500, '$/bins (450 kilos)', 'Región de O''Higgins', 661, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44655, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Segunda', 16, 255000, 260000, 257500, '$/bins (450 kilos)', 'Región de O''Higgins', 572, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44649, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 18, 330000, 340000, 335000, '$/bins (450 kilos)', 'Región de O''Higgins', 744, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44622, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 16, 410000, 420000, 415000, '$/bins (450 kilos)', 'Región de O''Higgins', 922, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44316, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 20, 255000, 260000, 257500, '$/bins (450 kilos)', 'Región de O''Higgins', 572, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44316, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 20, 225000, 230000, 227500, '$/bins (450 kilos)', 'Región de O''Higgins', 506, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44627, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 16, 405000, 410000, 407500, '$/bins (450 kilos)', 'Región de O''Higgins', 906, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44634, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 10, 305000, 310000, 307500, '$/bins (450 kilos)', 'Región Metropolitana', 683, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44648, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 16, 335000, 340000, 337500, '$/bins (450 kilos)', 'Región de O''Higgins', 750, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44273, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 10, 255000, 260000, 257500, '$/bins (450 kilos)', 'Región de O''Higgins', 572, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44273, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 20, 225000, 230000, 227500, '$/bins (450 kilos)', 'Región de O''Higgins', 506, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44628, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 14, 400000, 410000, 405000, '$/bins (450 kilos)', 'Región de O''Higgins', 900, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44294, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 20, 225000, 230000, 227500, '$/bins (450 kilos)', 'Región Metropolitana', 506, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44294, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Segunda', 16, 195000, 200000, 197500, '$/bins (450 kilos)', 'Región Metropolitana', 439, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44635, 4, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 20, 300000, 310000, 305000, '$/bins (450 kilos)', 'Región Metropolitana', 678, 450),
)

$dateCols = @(4)  # column D (1-indexed) holds the date serial value

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $colNum = $c + 1
        $cell = $ws.Cells.Item($rowNum, $colNum)
        if ($dateCols -contains $colNum) {
            $cell.NumberFormat = 'YYYY-MM-DD HH:MM:SS'
        }
        $cell.Value = $rowVals[$c]
    }
}
